$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new application record as row 5: Application ID, Project ID,
# Applicant NRIC, Application Status, Flat Type, Date
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 1).Style = "Normal"

$ws.Cells.Item(5, 2).Value = 3
$ws.Cells.Item(5, 2).Style = "Normal"

$ws.Cells.Item(5, 3).Value = "S1234567A"
$ws.Cells.Item(5, 5).Value = "2-ROOM"
$ws.Cells.Item(5, 4).Value = "Successful"

$ws.Cells.Item(5, 6).Value = 45770.31946777778
